$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark additional checklist cells as "x" (consolidating progress into the report)
$ws.Range("H2").Value = "x"
$ws.Range("C5").Value = "x"
$ws.Range("C6").Value = "x"
$ws.Range("B8").Value = "x"
$ws.Range("B10").Value = "x"

# Restore the scrolled/selected view state saved with the workbook
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 261
$ws.Range("H4").Select()
